$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 2 data (the "database" record that was fixed)
$ws.Range("A2").Value = 1222
$ws.Range("B2").Value = "jawa orang"
$ws.Range("C2").Value = "perempuan"
$ws.Range("D2").Value = "12/12/2000"
$ws.Range("E2").Value = "solokan"
$ws.Range("F2").Value = 9128378
$ws.Range("G2").Value = 1212321
$ws.Range("H2").Value = "dani@gmail.com"
$ws.Range("I2").Value = "12/12/2020"
$ws.Range("J2").Value = "Tetap"

# Add hyperlink on the email cell
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:dani@gmail.com")
$ws.Range("H2").Style = "Hyperlink"

# Update selection to match the final saved state
$ws.Range("F10").Select()
